$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.951.64'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '2.331.33'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("E6").Value = '  -2.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.511'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.07%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -4.23%  '
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.45'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0796'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.48%  '
$ws.Range("D16").Value = '2.318.76'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '42.885.90'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("E21").Value = '  -2.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.70'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.61'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.84%  '
$ws.Range("E38").Value = '  -2.29%  '
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("E40").Value = '  -5.04%  '
$ws.Range("E41").Value = '  -3.77%  '
$ws.Range("E42").Value = '  -2.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.31%  '
$ws.Range("D44").Value = '2.014.41'
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").Value = '2.557.30'
$ws.Range("E51").Value = '  +1.02%  '
